# Added for Program module
$wb = $excel.ActiveWorkbook

# Remove the empty "Sheet2" tab that sat between "Login" and "Program"
$emptySheet = $wb.Worksheets.Item("Sheet2")
$emptySheet.Delete()

# Make "Program" the active sheet (it becomes index 2 / activeTab 1 now that Sheet2 is gone)
$program = $wb.Worksheets.Item("Program")
$program.Activate()

# Populate the new Program module rows.
# Values are entered in this particular order so the shared-string table
# ends up with the same ordering as the authored workbook.
$program.Range("A4").Value = 123
$program.Range("B4").Value = 1345

$program.Range("C5").Value = "Active"
$program.Range("B8").Value = "4232#"
$program.Range("B5").Value = "ZSY123 Desc"
$program.Range("A5").Value = "CZX123"
$program.Range("A6").Value = "DAZXS"
$program.Range("B7").Value = "ABC Description123"
$program.Range("A8").Value = 890364

# Bold the header row
$program.Range("A1:C1").Font.Bold = $true

# Leave the active cell on A8, matching the authored selection
$program.Range("A8").Select()
